$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.060.12"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.832.21"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  +0.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.43"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6174"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9992"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07452"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2923"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.04"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07664"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "1.834.43"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.997"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6732"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.83"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009167"
$ws.Range("E16").Value = "  -4.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.915"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "29.043.78"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "2.074.03"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.59"
$ws.Range("E20").Value = "  +5.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.69"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.206"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9984"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.96"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1407"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.500"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.88"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.498"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05611"
$ws.Range("E30").Value = "  +3.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.138"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.111"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.208"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.840"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7411"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.653"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01784"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").Value = "1.213.58"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.400"
$ws.Range("E41").Value = "  -3.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8968"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9978"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.41"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "1.974.18"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5080"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000120"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4062"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.142"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05813"
$ws.Range("E51").Value = "  +0.60%  "
